$wb = $excel.ActiveWorkbook

# The workbook has two sheets with identical data tables that both need
# updating: "展览" and "全部类型". ("演出" and "本地生活" only contain headers.)
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 2277
    $ws.Range("F6").Value = 847
    $ws.Range("F8").Value = 5843
}
